$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4646496666666666
$ws.Range("H2").Value = 1.393949
$ws.Range("I2").Value = 0.08919948595155881
$ws.Range("J2").Value = 0.08919948595155881
$ws.Range("M2").Value = 0.303146
$ws.Range("N2").Value = 0.909438
$ws.Range("O2").Value = 0.005142855213700541
$ws.Range("P2").Value = 0.005142855213700542
$ws.Range("Q2").Value = 0.1408566878513333
$ws.Range("R2").Value = 1.267710190662
$ws.Range("S2").Value = 0.0004587400413853824
$ws.Range("T2").Value = 0.0004587400413853824

$ws.Range("G3").Value = 0.4646496666666666
$ws.Range("H3").Value = 1.393949
$ws.Range("I3").Value = 0.08919948595155881
$ws.Range("J3").Value = 0.08919948595155881
$ws.Range("O3").Value = 0.2877784259203595
$ws.Range("P3").Value = 0.2877784259203595
$ws.Range("Q3").Value = 7.881908828042777
$ws.Range("R3").Value = 70.937179452385
$ws.Range("S3").Value = 0.02566968766004481
$ws.Range("T3").Value = 0.02566968766004481

$ws.Range("G4").Value = 0.4646496666666666
$ws.Range("H4").Value = 1.393949
$ws.Range("I4").Value = 0.08919948595155881
$ws.Range("J4").Value = 0.08919948595155881
$ws.Range("M4").Value = 41.67881
$ws.Range("N4").Value = 125.03643
$ws.Range("O4").Value = 0.7070787188659401
$ws.Range("P4").Value = 0.7070787188659401
$ws.Range("Q4").Value = 19.36604517356333
$ws.Range("R4").Value = 174.29440656207
$ws.Range("S4").Value = 0.06307105825012863
$ws.Range("T4").Value = 0.06307105825012863

$ws.Range("I5").Value = 0.60288855652226
$ws.Range("J5").Value = 0.60288855652226
$ws.Range("M5").Value = 0.303146
$ws.Range("N5").Value = 0.909438
$ws.Range("O5").Value = 0.005142855213700541
$ws.Range("P5").Value = 0.005142855213700542
$ws.Range("Q5").Value = 0.9520333476059999
$ws.Range("R5").Value = 8.568300128454
$ws.Range("S5").Value = 0.003100568556190898
$ws.Range("T5").Value = 0.003100568556190899

$ws.Range("I6").Value = 0.60288855652226
$ws.Range("J6").Value = 0.60288855652226
$ws.Range("O6").Value = 0.2877784259203595
$ws.Range("P6").Value = 0.2877784259203595
$ws.Range("S6").Value = 0.1734983198013736
$ws.Range("T6").Value = 0.1734983198013736

$ws.Range("I7").Value = 0.60288855652226
$ws.Range("J7").Value = 0.60288855652226
$ws.Range("M7").Value = 41.67881
$ws.Range("N7").Value = 125.03643
$ws.Range("O7").Value = 0.7070787188659401
$ws.Range("P7").Value = 0.7070787188659401
$ws.Range("Q7").Value = 130.89276127191
$ws.Range("R7").Value = 1178.03485144719
$ws.Range("S7").Value = 0.4262896681646955
$ws.Range("T7").Value = 0.4262896681646955

$ws.Range("G8").Value = 1.603946333333333
$ws.Range("H8").Value = 4.811839
$ws.Range("I8").Value = 0.3079119575261812
$ws.Range("J8").Value = 0.3079119575261813
$ws.Range("M8").Value = 0.303146
$ws.Range("N8").Value = 0.909438
$ws.Range("O8").Value = 0.005142855213700541
$ws.Range("P8").Value = 0.005142855213700542
$ws.Range("Q8").Value = 0.4862299151646666
$ws.Range("R8").Value = 4.376069236482
$ws.Range("S8").Value = 0.00158354661612426
$ws.Range("T8").Value = 0.001583546616124261

$ws.Range("G9").Value = 1.603946333333333
$ws.Range("H9").Value = 4.811839
$ws.Range("I9").Value = 0.3079119575261812
$ws.Range("J9").Value = 0.3079119575261813
$ws.Range("O9").Value = 0.2877784259203595
$ws.Range("P9").Value = 0.2877784259203595
$ws.Range("Q9").Value = 27.20793679913722
$ws.Range("R9").Value = 244.871431192235
$ws.Range("S9").Value = 0.08861041845894101
$ws.Range("T9").Value = 0.08861041845894103

$ws.Range("G10").Value = 1.603946333333333
$ws.Range("H10").Value = 4.811839
$ws.Range("I10").Value = 0.3079119575261812
$ws.Range("J10").Value = 0.3079119575261813
$ws.Range("M10").Value = 41.67881
$ws.Range("N10").Value = 125.03643
$ws.Range("O10").Value = 0.7070787188659401
$ws.Range("P10").Value = 0.7070787188659401
$ws.Range("Q10").Value = 66.85057447719666
$ws.Range("R10").Value = 601.65517029477
$ws.Range("S10").Value = 0.217717992451116
$ws.Range("T10").Value = 0.217717992451116
